# interfaz_manager.xlsx maintenance pass
# - rename the sheet to match the file name
# - move the active selection to I25 (where the user left off)
# - re-apply center/center alignment to the "tipo_dato" (column D) cells
#   that were still carrying the old, redundant fill-flagged style so
#   Excel can fold it into the identical fill-less style already in use
#   elsewhere in the column

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "interfaz_manager"

$xlCenter = -4108

$dRows = @(9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,31,39,42,43,48)
foreach ($r in $dRows) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.HorizontalAlignment = $xlCenter
    $cell.VerticalAlignment = $xlCenter
}

$ws.Range("I25").Select() | Out-Null
